$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Widen the metric/value columns
# (ColumnWidth -> stored XML width has a constant +5/6 offset in this engine,
#  so back the requested value off by 5/6 to land on the exact target width)
$ws.Columns.Item(1).ColumnWidth = 38 - 5/6
$ws.Columns.Item(2).ColumnWidth = 11 - 5/6

# Row 2: Total Appointments count
$ws.Range("B2").Value = 11

# Row 3: Dr. Shreyansh count
$ws.Range("B3").Value = 5

# Row 4: was "Appointments at Jubliee Hills" -> now Dr. Naveen (Doctor Distribution)
$ws.Range("A4").Value = "Appointments with Dr. Naveen"
$ws.Range("B4").Value = 2
$ws.Range("C4").Value = "Doctor Distribution"

# Row 5: was "Patients with Aetna" -> now Dr. Aish (Doctor Distribution)
$ws.Range("A5").Value = "Appointments with Dr. Aish"
$ws.Range("B5").Value = 2
$ws.Range("C5").Value = "Doctor Distribution"

# Row 6 (new): Dr. Naresh
$ws.Range("A6").Value = "Appointments with Dr. Naresh"
$ws.Range("B6").Value = 2
$ws.Range("C6").Value = "Doctor Distribution"

# Row 7 (new): Banjara Hills
$ws.Range("A7").Value = "Appointments at Banjara Hills"
$ws.Range("B7").Value = 5
$ws.Range("C7").Value = "Location Distribution"

# Row 8 (new): Gachibowli
$ws.Range("A8").Value = "Appointments at Gachibowli"
$ws.Range("B8").Value = 3
$ws.Range("C8").Value = "Location Distribution"

# Row 9 (new): Jubliee Hills
$ws.Range("A9").Value = "Appointments at Jubliee Hills"
$ws.Range("B9").Value = 3
$ws.Range("C9").Value = "Location Distribution"

# Row 10 (new): Patients with Aetna
$ws.Range("A10").Value = "Patients with Aetna"
$ws.Range("B10").Value = 4
$ws.Range("C10").Value = "Insurance Analysis"

# Row 11 (new): Blue Cross Blue Shield
$ws.Range("A11").Value = "Patients with Blue Cross Blue Shield"
$ws.Range("B11").Value = 3
$ws.Range("C11").Value = "Insurance Analysis"

# Row 12 (new): Cigna
$ws.Range("A12").Value = "Patients with Cigna"
$ws.Range("B12").Value = 2
$ws.Range("C12").Value = "Insurance Analysis"

# Row 13 (new): aetna (lowercase)
$ws.Range("A13").Value = "Patients with aetna"
$ws.Range("B13").Value = 1
$ws.Range("C13").Value = "Insurance Analysis"

# Row 14 (new): cigna (lowercase)
$ws.Range("A14").Value = "Patients with cigna"
$ws.Range("B14").Value = 1
$ws.Range("C14").Value = "Insurance Analysis"

# Row 15: Estimated Revenue (was row 6), updated amount
$ws.Range("A15").Value = "Estimated Revenue"
$ws.Range("B15").Value = "'$1,020.00"
$ws.Range("C15").Value = "Financial"

# Copy the style used by data rows (row 2, e.g. A2) down to all the newly added rows
$ws.Range("A2:C2").Copy() | Out-Null
$ws.Range("A3:C15").PasteSpecial(-4122) | Out-Null
